$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "(0.02)"
$ws.Range("C4").Value = "(0.14)"
$ws.Range("D4").Value = "(0.09)"
$ws.Range("E4").Value = "(0.03)"
$ws.Range("F4").Value = "(0.09)"
$ws.Range("G4").Value = "(0.13)"

$ws.Range("B6").Value = "(0.01)"
$ws.Range("C6").Value = "(0.13)"
$ws.Range("D6").Value = "(0.08)"
$ws.Range("E6").Value = "(0.06)"
$ws.Range("F6").Value = "(0.23)"
$ws.Range("G6").Value = "(0.33)"
